$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 136; this shifts the existing rows 136-207
# down to 137-208, preserving all their data untouched.
$ws.Rows.Item(136).Insert()

# Populate the newly inserted row 136 with the new weekly price record.
$ws.Range("A136").Value = 6
$ws.Range("B136").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C136").Value = "Metropolitana"
$ws.Range("D136").Value = 44466
$ws.Range("E136").Value = 13
$ws.Range("F136").Value = 100112032
$ws.Range("G136").Value = "Zapallo italiano"
$ws.Range("H136").Value = "Sin especificar"
$ws.Range("I136").Value = "Primera"
$ws.Range("J136").Value = 490
$ws.Range("K136").Value = 12000
$ws.Range("L136").Value = 14000
$ws.Range("M136").Value = 12939
$ws.Range("N136").Value = "$/caja 50 unidades"
$ws.Range("O136").Value = "Región de Arica y Parinacota"
$ws.Range("P136").Value = 259
$ws.Range("Q136").Value = 50
$ws.Range("R136").Value = "Hortaliza"
